$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" row (A8 = "Date") with the new timestamp value.
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# Update the "Case Sensitive" row (A20 = "Case Sensitive") with the text
# value "true". A plain assignment of the literal "true" would be
# auto-coerced into an Excel boolean TRUE, so instead we build it via a
# formula that evaluates to the text string "true" and then convert that
# formula result to a static value in place (Copy + PasteSpecial values),
# which preserves the existing cell style and stores it as a normal
# shared-string text cell instead of a boolean.
$r = $ws.Range("B20")
$r.Formula = '="true"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0
